$wb = $excel.ActiveWorkbook
$customers = $wb.Worksheets.Item("customers")

# Add the new "addresses" worksheet right after "customers"
$addresses = $wb.Worksheets.Add($null, $customers)
$addresses.Name = "addresses"

# Header row
$addresses.Range("A1").Value = "customerId"
$addresses.Range("B1").Value = "number"
$addresses.Range("C1").Value = "street"
$addresses.Range("D1").Value = "city"
$addresses.Range("E1").Value = "country"

# Data row
$addresses.Range("A2").Value = 1
$addresses.Range("B2").Value = 42
$addresses.Range("C2").Value = "expectedStreet"
$addresses.Range("D2").Value = "expectedCIty"
$addresses.Range("E2").Value = "expectedCountry"

$addresses.Range("E3").Select() | Out-Null
